# Generate Report for Handback
# Updates the handback-status report with refreshed timestamps / status
# for the 5ecf8e06-7634-4602-9ba4-8caf09518d61 row across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-12 06:19:21"
$wsOverview.Range("G5").Value = "2016-08-12 06:19:21"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-12 06:19:13"
$wsZhCn.Range("H5").Value = "2016-08-12 06:19:13"
$wsZhCn.Range("K3").Value = "2016-08-12 06:19:41"
$wsZhCn.Range("K5").Value = "2016-08-12 06:19:41"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-12 06:19:21"
$wsDeDe.Range("H5").Value = "2016-08-12 06:19:21"
$wsDeDe.Range("K3").Value = "2016-08-12 06:19:52"
$wsDeDe.Range("K5").Value = "2016-08-12 06:19:52"
